# Apply the FLORIDA_2017 sheet1 corrections:
#  - rename header row to snake_case field names
#  - title-case Spanish connector words (de/del/el/los/las/la/y) in
#    mx_state / mx_municipality text
#  - two floating point value corrections (recomputed percentages)
#  - drop the trailing footnote rows (1721:1726)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @('A1', 'mx_state'),
    @('B1', 'mx_municipality'),
    @('C1', 'n_matriculas'),
    @('D1', 'pct_matriculas'),
    @('B8', 'Pabellón De Arteaga'),
    @('B9', 'Rincón De Romos'),
    @('B10', 'San Francisco De Los Romo'),
    @('B11', 'San José De Gracia'),
    @('B38', 'Amatenango De La Frontera'),
    @('B41', 'Bejucal De Ocampo'),
    @('B51', 'Chiapa De Corzo'),
    @('B58', 'Comitán De Domínguez'),
    @('B85', 'Marqués De Comillas'),
    @('B86', 'Mazapa De Madero'),
    @('B91', 'Ocozocoautla De Espinosa'),
    @('B103', 'Salto De Agua'),
    @('B104', 'San Cristóbal De Las Casas'),
    @('B146', 'Guadalupe Y Calvo'),
    @('B147', 'Hidalgo Del Parral'),
    @('B156', 'San Francisco De Borja'),
    @('B178', 'San Juan De Sabinas'),
    @('D184', 0.009670211126521808),
    @('B192', 'Villa De Álvarez'),
    @('A194', 'Ciudad De México'),
    @('B198', 'Cuajimalpa De Morelos'),
    @('B212', 'Coneto De Comonfort'),
    @('B226', 'Nombre De Dios'),
    @('B230', 'Pánuco De Coronado'),
    @('B236', 'San Juan De Guadalupe'),
    @('B237', 'San Juan Del Río'),
    @('B238', 'San Luis Del Cordero'),
    @('A247', 'Estado De México'),
    @('B247', 'Acambay De Ruíz Castañeda'),
    @('B250', 'Almoloya De Alquisiras'),
    @('B251', 'Almoloya De Juárez'),
    @('B258', 'Atizapán De Zaragoza'),
    @('B263', 'Chapa De Mota'),
    @('B268', 'Coacalco De Berriozábal'),
    @('B275', 'Ecatepec De Morelos'),
    @('B282', 'Ixtapan De La Sal'),
    @('B295', 'Naucalpan De Juárez'),
    @('B307', 'San Felipe Del Progreso'),
    @('B309', 'San Simón De Guerrero'),
    @('B311', 'Soyaniquilpan De Juárez'),
    @('B319', 'Tenango Del Aire'),
    @('B320', 'Tenango Del Valle'),
    @('B332', 'Tlalnepantla De Baz'),
    @('B338', 'Valle De Bravo'),
    @('B339', 'Valle De Chalco Solidaridad'),
    @('B340', 'Villa De Allende'),
    @('B341', 'Villa Del Carbón'),
    @('B353', 'San Miguel De Allende'),
    @('B354', 'Apaseo El Alto'),
    @('B355', 'Apaseo El Grande'),
    @('B363', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B367', 'Jaral Del Progreso'),
    @('B375', 'Purísima Del Rincón'),
    @('B379', 'San Diego De La Unión'),
    @('B381', 'San Francisco Del Rincón'),
    @('B383', 'San Luis De La Paz'),
    @('B384', 'Santa Cruz De Juventino Rosas'),
    @('B386', 'Silao De La Victoria'),
    @('B391', 'Valle De Santiago'),
    @('B397', 'Acapulco De Juárez'),
    @('B400', 'Ajuchitlán Del Progreso'),
    @('B401', 'Alcozauca De Guerrero'),
    @('B405', 'Atenango Del Río'),
    @('B407', 'Atoyac De Álvarez'),
    @('B408', 'Ayutla De Los Libres'),
    @('B411', 'Buenavista De Cuéllar'),
    @('B412', 'Chilapa De Álvarez'),
    @('B413', 'Chilpancingo De Los Bravo'),
    @('B414', 'Coahuayutla De José María Izazaga'),
    @('B419', 'Coyuca De Benítez'),
    @('B420', 'Coyuca De Catalán'),
    @('B424', 'Cuetzala Del Progreso'),
    @('B425', 'Cutzamala De Pinzón'),
    @('B431', 'Huitzuco De Los Figueroa'),
    @('B432', 'Iguala De La Independencia'),
    @('B434', 'Ixcateopan De Cuauhtémoc'),
    @('B435', 'Zihuatanejo De Azueta'),
    @('B437', 'La Unión De Isidoro Montes De Oca'),
    @('B440', 'Mártir De Cuilapan'),
    @('B453', 'Taxco De Alarcón'),
    @('B455', 'Técpan De Galeana'),
    @('B457', 'Tepecoacuilco De Trujano'),
    @('B459', 'Tixtla De Guerrero'),
    @('B463', 'Tlalixtaquilla De Maldonado'),
    @('B464', 'Tlapa De Comonfort'),
    @('D472', 0.09597010325165664),
    @('B476', 'Agua Blanca De Iturbide'),
    @('B482', 'Atotonilco De Tula'),
    @('B483', 'Atotonilco El Grande'),
    @('B489', 'Cuautepec De Hinojosa'),
    @('B495', 'Huasca De Ocampo'),
    @('B499', 'Huejutla De Reyes'),
    @('B502', 'Jacala De Ledezma'),
    @('B508', 'Mineral Del Chico'),
    @('B509', 'Mineral Del Monte'),
    @('B510', 'Mixquiahuala De Juárez'),
    @('B511', 'Molango De Escamilla'),
    @('B513', 'Nopala De Villagrán'),
    @('B514', 'Omitlán De Juárez'),
    @('B515', 'Pachuca De Soto'),
    @('B518', 'Progreso De Obregón'),
    @('B524', 'Santiago De Anaya'),
    @('B525', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B529', 'Tenango De Doria'),
    @('B531', 'Tepehuacán De Guerrero'),
    @('B532', 'Tepeji Del Río De Ocampo'),
    @('B535', 'Tezontepec De Aldama'),
    @('B542', 'Tula De Allende'),
    @('B543', 'Tulancingo De Bravo'),
    @('B544', 'Villa De Tezontepec'),
    @('B547', 'Zacualtipán De Ángeles'),
    @('B551', 'Acatlán De Juárez'),
    @('B552', 'Ahualulco De Mercado'),
    @('B557', 'Atotonilco El Alto'),
    @('B558', 'Autlán De Navarro'),
    @('B563', 'Cañadas De Obregón'),
    @('B568', 'Concepción De Buenos Aires'),
    @('B576', 'Encarnación De Díaz'),
    @('B582', 'Huejuquilla El Alto'),
    @('B583', 'Ixtlahuacán Del Río'),
    @('B587', 'Jilotlán De Los Dolores'),
    @('B593', 'La Manzanilla De La Paz'),
    @('B594', 'Lagos De Moreno'),
    @('B600', 'Ojuelos De Jalisco'),
    @('B605', 'San Diego De Alejandría'),
    @('B606', 'San Juan De Los Lagos'),
    @('B608', 'San Sebastián Del Oeste'),
    @('B611', 'Tamazula De Gordiano'),
    @('B613', 'Techaluta De Montenegro'),
    @('B617', 'Teocuitatlán De Corona'),
    @('B618', 'Tepatitlán De Morelos'),
    @('B620', 'Tlajomulco De Zúñiga'),
    @('B628', 'Unión De San Antonio'),
    @('B629', 'Unión De Tula'),
    @('B630', 'Valle De Juárez'),
    @('B635', 'Yahualica De González Gallo'),
    @('B638', 'Zapotitlán De Vadillo'),
    @('B639', 'Zapotlán Del Rey'),
    @('B640', 'Zapotlán El Grande'),
    @('B666', 'Coalcomán De Vázquez Pallares'),
    @('B668', 'Cojumatlán De Régules'),
    @('B734', 'Tiquicheo De Nicolás Romero'),
    @('B760', 'Coatlán Del Río'),
    @('B768', 'Jonacatepec De Leandro Valle'),
    @('B772', 'Puente De Ixtla'),
    @('B777', 'Tetela Del Volcán'),
    @('B778', 'Tlaltizapán De Zapata'),
    @('B789', 'Bahía De Banderas'),
    @('B792', 'Ixtlán Del Río'),
    @('B799', 'Santa María Del Oro'),
    @('B812', 'Ciénega De Flores'),
    @('B824', 'Mier Y Noriega'),
    @('B829', 'San Nicolás De Los Garza'),
    @('B834', 'Acatlán De Pérez Figueroa'),
    @('B840', 'Chalcatongo De Hidalgo'),
    @('B842', 'Coicoyán De Las Flores'),
    @('B843', 'Constancia Del Rosario'),
    @('B845', 'Fresnillo De Trujano'),
    @('B846', 'Guadalupe De Ramírez'),
    @('B848', 'Guevea De Humboldt'),
    @('B849', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B850', 'Heroica Ciudad De Huajuapan De León'),
    @('B851', 'Heroica Ciudad De Tlaxiaco'),
    @('B852', 'Huautla De Jiménez'),
    @('B853', 'Ixtlán De Juárez'),
    @('B854', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B863', 'Mariscala De Juárez'),
    @('B865', 'Mazatlán Villa De Flores'),
    @('B866', 'Miahuatlán De Porfirio Díaz'),
    @('B867', 'Mixistlán De La Reforma'),
    @('B870', 'Nejapa De Madero'),
    @('B871', 'Oaxaca De Juárez'),
    @('B872', 'Ocotlán De Morelos'),
    @('B873', 'Pinotepa De Don Luis'),
    @('B875', 'Putla Villa De Guerrero'),
    @('B877', 'Rojas De Cuauhtémoc'),
    @('B889', 'San Antonino El Alto'),
    @('B907', 'San Dionisio Del Mar'),
    @('B912', 'San Francisco Del Mar'),
    @('B929', 'San Juan Bautista Lo De Soto'),
    @('B937', 'San Juan De Los Cués'),
    @('B938', 'San Juan Del Río'),
    @('B982', 'San Miguel Del Puerto'),
    @('B983', 'San Miguel Del Río'),
    @('B984', 'San Miguel El Grande'),
    @('B997', 'San Pablo Villa De Mitla'),
    @('B1001', 'San Pedro El Alto'),
    @('B1012', 'San Pedro Y San Pablo Ayutla'),
    @('B1013', 'San Pedro Y San Pablo Teposcolula'),
    @('B1022', 'Santa Ana Del Valle'),
    @('B1037', 'Santa Cruz Tacache De Mina'),
    @('B1043', 'Santa Inés De Zaragoza'),
    @('B1044', 'Santa Inés Del Monte'),
    @('B1046', 'Santa Lucía Del Camino'),
    @('B1056', 'Santa María Del Rosario'),
    @('B1062', 'Santa María Jalapa Del Marqués'),
    @('B1112', 'Santo Domingo De Morelos'),
    @('B1128', 'Tamazulápam Del Espíritu Santo'),
    @('B1129', 'Tataltepec De Valdés'),
    @('B1130', 'Teotitlán De Flores Magón'),
    @('B1131', 'Teotitlán Del Valle'),
    @('B1133', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'),
    @('B1134', 'Tlacolula De Matamoros'),
    @('B1135', 'Totontepec Villa De Morelos'),
    @('B1139', 'Villa De Etla'),
    @('B1140', 'Villa De Tututepec'),
    @('B1141', 'Villa De Zaachila'),
    @('B1144', 'Villa Sola De Vega'),
    @('B1145', 'Zapotitlán Del Río'),
    @('B1148', 'Zimatlán De Álvarez'),
    @('B1166', 'Ayotoxco De Guerrero'),
    @('B1178', 'Chila De La Sal'),
    @('B1187', 'Cuayuca De Andrade'),
    @('B1188', 'Cuetzalan Del Progreso'),
    @('B1202', 'Huehuetlán El Chico'),
    @('B1203', 'Huehuetlán El Grande'),
    @('B1208', 'Huitzilan De Serdán'),
    @('B1210', 'Ixcamilpa De Guerrero'),
    @('B1213', 'Izúcar De Matamoros'),
    @('B1222', 'Los Reyes De Juárez'),
    @('B1229', 'Palmar De Bravo'),
    @('B1242', 'San Nicolás De Los Ranchos'),
    @('B1245', 'San Salvador El Seco'),
    @('B1246', 'San Salvador El Verde'),
    @('B1255', 'Tepanco De López'),
    @('B1256', 'Tepango De Rodríguez'),
    @('B1260', 'Tepexi De Rodríguez'),
    @('B1262', 'Tetela De Ocampo'),
    @('B1267', 'Tlacotepec De Benito Juárez'),
    @('B1277', 'Tuzamapan De Galeana'),
    @('B1280', 'Xayacatlán De Bravo'),
    @('B1283', 'Xochitlán De Vicente Suárez'),
    @('B1295', 'Amealco De Bonfil'),
    @('B1297', 'Cadereyta De Montes'),
    @('B1303', 'Jalpan De Serra'),
    @('B1304', 'Landa De Matamoros'),
    @('B1307', 'Pinal De Amoles'),
    @('B1310', 'San Juan Del Río'),
    @('B1322', 'Armadillo De Los Infante'),
    @('B1323', 'Axtla De Terrazas'),
    @('B1328', 'Ciudad Del Maíz'),
    @('B1338', 'Mexquitic De Carmona'),
    @('B1343', 'San Ciro De Acosta'),
    @('B1349', 'Santa María Del Río'),
    @('B1351', 'Soledad De Graciano Sánchez'),
    @('B1359', 'Tanquián De Escobedo'),
    @('B1363', 'Villa De Arista'),
    @('B1364', 'Villa De Arriaga'),
    @('B1365', 'Villa De Guadalupe'),
    @('B1366', 'Villa De La Paz'),
    @('B1367', 'Villa De Ramos'),
    @('B1368', 'Villa De Reyes'),
    @('B1403', 'Nacozari De García'),
    @('B1419', 'Jalpa De Méndez'),
    @('B1460', 'Soto La Marina'),
    @('B1467', 'Acuamanala De Miguel Hidalgo'),
    @('B1468', 'Apetatitlán De Antonio Carvajal'),
    @('B1472', 'Contla De Juan Cuamatzi'),
    @('B1477', 'Ixtacuixtla De Mariano Matamoros'),
    @('B1479', 'Muñoz De Domingo Arenas'),
    @('B1480', 'Nanacamilpa De Mariano Arista'),
    @('B1483', 'Papalotla De Xicohténcatl'),
    @('B1484', 'San Pablo Del Monte'),
    @('B1500', 'Alto Lucero De Gutiérrez Barrios'),
    @('B1504', 'Amatlán De Los Reyes'),
    @('B1514', 'Boca Del Río'),
    @('B1515', 'Camarón De Tejeda'),
    @('B1518', 'Castillo De Teayo'),
    @('B1520', 'Cazones De Herrera'),
    @('B1536', 'Cosamaloapan De Carpio'),
    @('B1537', 'Cosautlán De Carvajal'),
    @('B1553', 'Hueyapan De Ocampo'),
    @('B1554', 'Ignacio De La Llave'),
    @('B1558', 'Ixhuatlán De Madero'),
    @('B1559', 'Ixhuatlán Del Café'),
    @('B1570', 'Juchique De Ferrer'),
    @('B1573', 'Landero Y Coss'),
    @('B1575', 'Las Vigas De Ramírez'),
    @('B1576', 'Lerdo De Tejada'),
    @('B1579', 'Martínez De La Torre'),
    @('B1581', 'Medellín De Bravo'),
    @('B1585', 'Mixtla De Altamirano'),
    @('B1598', 'Paso De Ovejas'),
    @('B1599', 'Paso Del Macho'),
    @('B1602', 'Poza Rica De Hidalgo'),
    @('B1610', 'Sayula De Alemán'),
    @('B1613', 'Soledad De Doblado'),
    @('B1635', 'Tlacotepec De Mejía'),
    @('B1648', 'Vega De Alatorre'),
    @('B1658', 'Zontecomatlán De López Y Fuentes'),
    @('B1680', 'Cañitas De Felipe Pescador'),
    @('B1681', 'Concepción Del Oro'),
    @('B1683', 'El Plateado De Joaquín Amaro'),
    @('B1692', 'Jiménez Del Teul'),
    @('B1699', 'Nochistlán De Mejía'),
    @('B1700', 'Noria De Ángeles'),
    @('B1709', 'Tlaltenango De Sánchez Román'),
    @('B1713', 'Villa De Cos')
)

foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# Remove the sample-size / source / author / date footnote rows that
# trailed the data table; this also shrinks the sheet dimension from
# A1:D1726 down to A1:D1720.
$ws.Rows("1721:1726").Delete()
